# GradedExercise - "Good naming" criterion (row 20) is now fully awarded:
# mark D20/E20 with the "Neutral" highlight style (matching the other
# scored sub-criteria rows 17-19) and bump the achieved Value (G20) from
# 0 up to its max of 2. Downstream subtotal / total formulas recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the "Good naming" row like its already-graded neighbours.
$ws.Range("D20:E20").Style = "Neutral"

# Award full points for this criterion.
$ws.Range("G20").Value = 2

# Move the on-screen selection / viewport the way the author left it.
$ws.Activate()
$ws.Range("L20").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1

# Make sure every dependent formula (subtotals, bonus/malus, grade) is fresh.
$excel.CalculateFull()
